# Phoenix roster: re-run of the totals/prediction pipeline swapped the
# per-player detail (No., Player, Pos, Ht, Wt, Birth Date, Nationality,
# Exp, College, bbref url) between a handful of row pairs, while the
# leading index column (A) stayed put. Also:
#   - "Ish Wainright (TW)" lost its "(TW)" two-way suffix.
#   - Kevin Durant / Darius Bazley (rows 16-17) picked up jersey numbers
#     (35 / 55) that were previously blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 7 <-> Row 8 (Bismack Biyombo <-> Ish Wainright) ----
$ws.Range("B7").Value2 = 12
$ws.Range("C7").Value2 = "Ish Wainright"
$ws.Range("D7").Value2 = "SF"
$ws.Range("E7").Value2 = "6-6"
$ws.Range("F7").Value2 = 250
$ws.Range("G7").Value2 = "September 12, 1994"
$ws.Range("H7").Value2 = "us"
$ws.Range("I7").Value2 = 1
$ws.Range("J7").Value2 = "Baylor"
$ws.Range("K7").Value2 = "https://www.basketball-reference.com/players/w/wainris01.html"

$ws.Range("B8").Value2 = 18
$ws.Range("C8").Value2 = "Bismack Biyombo"
$ws.Range("D8").Value2 = "C"
$ws.Range("E8").Value2 = "6-8"
$ws.Range("F8").Value2 = 255
$ws.Range("G8").Value2 = "August 28, 1992"
$ws.Range("H8").Value2 = "cd"
$ws.Range("I8").Value2 = 11
$ws.Range("J8").ClearContents()
$ws.Range("K8").Value2 = "https://www.basketball-reference.com/players/b/biyombi01.html"

# ---- Row 11 <-> Row 12 (Landry Shamet <-> Cameron Payne) ----
$ws.Range("B11").Value2 = 15
$ws.Range("C11").Value2 = "Cameron Payne"
$ws.Range("D11").Value2 = "PG"
$ws.Range("E11").Value2 = "6-1"
$ws.Range("F11").Value2 = 183
$ws.Range("G11").Value2 = "August 8, 1994"
$ws.Range("I11").Value2 = 7
$ws.Range("J11").Value2 = "Murray State"
$ws.Range("K11").Value2 = "https://www.basketball-reference.com/players/p/payneca01.html"

$ws.Range("B12").Value2 = 14
$ws.Range("C12").Value2 = "Landry Shamet"
$ws.Range("D12").Value2 = "SG"
$ws.Range("E12").Value2 = "6-4"
$ws.Range("F12").Value2 = 190
$ws.Range("G12").Value2 = "March 13, 1997"
$ws.Range("I12").Value2 = 4
$ws.Range("J12").Value2 = "Wichita State"
$ws.Range("K12").Value2 = "https://www.basketball-reference.com/players/s/shamela01.html"

# ---- Row 14 <-> Row 15 (T.J. Warren <-> Terrence Ross) ----
$ws.Range("B14").Value2 = 8
$ws.Range("C14").Value2 = "Terrence Ross"
$ws.Range("D14").Value2 = "SG"
$ws.Range("E14").Value2 = "6-6"
$ws.Range("F14").Value2 = 206
$ws.Range("G14").Value2 = "February 5, 1991"
$ws.Range("I14").Value2 = 10
$ws.Range("J14").Value2 = "Washington"
$ws.Range("K14").Value2 = "https://www.basketball-reference.com/players/r/rosste01.html"

$ws.Range("B15").Value2 = 21
$ws.Range("C15").Value2 = "T.J. Warren"
$ws.Range("D15").Value2 = "SF"
$ws.Range("E15").Value2 = "6-8"
$ws.Range("F15").Value2 = 220
$ws.Range("G15").Value2 = "September 5, 1993"
$ws.Range("I15").Value2 = 7
$ws.Range("J15").Value2 = "NC State"
$ws.Range("K15").Value2 = "https://www.basketball-reference.com/players/w/warretj01.html"

# ---- Row 16 <-> Row 17 (Darius Bazley <-> Kevin Durant) ----
$ws.Range("B16").Value2 = 35
$ws.Range("C16").Value2 = "Kevin Durant"
$ws.Range("D16").Value2 = "SF"
$ws.Range("E16").Value2 = "6-10"
$ws.Range("F16").Value2 = 240
$ws.Range("G16").Value2 = "September 29, 1988"
$ws.Range("I16").Value2 = 14
$ws.Range("J16").Value2 = "Texas"
$ws.Range("K16").Value2 = "https://www.basketball-reference.com/players/d/duranke01.html"

$ws.Range("B17").Value2 = 55
$ws.Range("C17").Value2 = "Darius Bazley"
$ws.Range("D17").Value2 = "PF"
$ws.Range("E17").Value2 = "6-8"
$ws.Range("F17").Value2 = 208
$ws.Range("G17").Value2 = "June 12, 2000"
$ws.Range("I17").Value2 = 3
$ws.Range("J17").ClearContents()
$ws.Range("K17").Value2 = "https://www.basketball-reference.com/players/b/bazleda01.html"
